{"js": "// Append \" (Changed main)\" to the end of the paragraph \"This is a\n// Microsoft word document.\" as three separate runs:\n//   1) \" (\"\n//   2) \"Changed main\"\n//   3) \")\"\n//\n// insertText() would merge the new characters into the existing run\n// (same formatting), so we use insertOoxml() with a Flat-OPC fragment\n// that explicitly defines three <w:r> elements \u2014 this preserves run\n// boundaries exactly as required.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text === \"This is a Microsoft word document.\"\n);\nif (!target) {\n  throw new Error(\"Could not find paragraph 'This is a Microsoft word document.'\");\n}\n\nconst endRange = target.getRange(\"End\");\n\nconst flatOpcXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>Changed main</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nendRange.insertOoxml(flatOpcXml, \"End\");\nawait context.sync();\n", "ps1": "# Append \" (Changed main)\" to the end of the paragraph \"This is a\n# Microsoft word document.\" as three separate runs:\n#   1) \" (\"\n#   2) \"Changed main\"\n#   3) \")\"\n#\n# Range.InsertAfter()/.Text would merge the new characters into the\n# existing run (same formatting), so instead we build a Flat-OPC OOXML\n# fragment with three explicit <w:r> elements and insert it with\n# InsertXML(..., \"End\"), which appends at the end of the target range\n# without disturbing (replacing) the run(s) already there.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"This is a Microsoft word document.\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'This is a Microsoft word document.'\"\n}\n\n$r = $target.Range\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n       '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n       '<pkg:xmlData>' +\n       '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n       '<w:body>' +\n       '<w:p>' +\n       '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n       '<w:r><w:t>Changed main</w:t></w:r>' +\n       '<w:r><w:t>)</w:t></w:r>' +\n       '</w:p>' +\n       '</w:body>' +\n       '</w:document>' +\n       '</pkg:xmlData>' +\n       '</pkg:part>' +\n       '</pkg:package>'\n\n$r.InsertXML($xml, \"End\")\n"}
